$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '57.777.15'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.40%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.126.45'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  +0.07%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '532.10'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.90%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '138.59'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("E7").Value = '  +0.02%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '3.125.96'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("E9").Value = '  +4.53%  '
$ws.Range("E10").Value = '  +1.21%  '
$ws.Range("E11").Value = '  +0.14%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.416'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +4.58%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '3.664.83'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("E14").Value = '  +1.38%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '25.51'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("E16").Value = '  +0.04%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '57.946.65'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.51%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.126.44'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.48%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.02'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("E20").Value = '  +0.68%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '8.11'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.35%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '360.40'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +3.15%  '
$ws.Range("E23").Value = '  -0.11%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '69.02'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("E25").Value = '  +0.09%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.167'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("E28").Value = '  -4.24%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.29'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("E30").Value = '  -0.17%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.08'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.13%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '21.38'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.40%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.12'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.79%  '
$ws.Range("E34").Value = '  -2.62%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '158.28'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("E36").Value = '  -1.43%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '25.86'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("E38").Value = '  +1.77%  '
$ws.Range("E39").Value = '  +2.44%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0673'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.92%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.496.51'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +6.42%  '
$ws.Range("E42").Value = '  -0.28%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.99'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -4.79%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '37.78'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +3.37%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.170.47'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("E47").Value = '  -0.38%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.987'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.97%  '
$ws.Range("E49").Value = '  +0.59%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '19.76'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.07%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.742'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.84%  '
